$d = $word.ActiveDocument

function Replace-ParaXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $rng = $p.Range.Duplicate
    $rng.Collapse(1)
    $pkg = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

function Insert-ParasAfter($paraIndex, $innerXml) {
    # collapses to a mid-paragraph point so the new paragraph(s) land right after
    # paragraph $paraIndex without consuming / replacing it.
    $p = $d.Paragraphs($paraIndex)
    $s = $p.Range.Start
    $e = $p.Range.End
    $mid = [Math]::Floor(($s + $e) / 2)
    if ($mid -eq $e) { $mid = $s }
    $rng = $d.Range($mid, $mid)
    $pkg = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# 1) Insert two new paragraphs ("Srodowisko programistyczne:" / "Visual Studio Code 1.68.0")
#    right after paragraph 22 ("3. Przebieg realizacji") and before "Pliki skladowe:".
$newParas = '<w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr><w:t>Środowisko programistyczne</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Visual Studio </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr><w:t>Code</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> 1.68.0</w:t></w:r></w:p>'
Insert-ParasAfter 22 $newParas

# From here on, everything from the old paragraph 23 onward is shifted down by 2.
$OFF = 2

# 2) Add the _GoBack bookmark to the end of the "Pliki skladowe:" paragraph
#    (old index 23 -> now 25).
$p23 = '<w:p><w:pPr><w:rPr><w:bCs/><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr><w:t>Pliki składowe:</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
Replace-ParaXml (23 + $OFF) $p23

# 3) Merge "functions_" + "menu.h" runs into a single "functions_menu.h" run
#    (old index 25 -> now 27).
$p25 = '<w:p><w:pPr><w:pStyle w:val="Bezodstpw"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="en-GB"/></w:rPr><w:t>H:</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>functions_game.h</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>functions_menu.h</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>functions_</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>system.h</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
Replace-ParaXml (25 + $OFF) $p25

# 4) Add <w:lastRenderedPageBreak/> before the "TXT:" run
#    (old index 27 -> now 29).
$p27 = '<w:p><w:pPr><w:pStyle w:val="Bezodstpw"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="en-GB"/></w:rPr><w:lastRenderedPageBreak/><w:t>TXT:</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> Bomb.txt, intro.txt, intro2.txt, logo.txt</w:t></w:r></w:p>'
Replace-ParaXml (27 + $OFF) $p27

# 5) Add <w:rPr><w:lang w:val="en-GB"/></w:rPr> to the two empty paragraphs
#    (old indices 29, 30 -> now 31, 32).
$pEmpty = '<w:p><w:pPr><w:pStyle w:val="Bezodstpw"/><w:ind w:firstLine="360"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p>'
Replace-ParaXml (29 + $OFF) $pEmpty
Replace-ParaXml (30 + $OFF) $pEmpty

# 6) Remove <w:lastRenderedPageBreak/> before "Realizacja projekt"
#    (old index 31 -> now 33).
$p31 = '<w:p><w:pPr><w:pStyle w:val="Bezodstpw"/><w:ind w:firstLine="360"/></w:pPr><w:r><w:t>Realizacja projekt</w:t></w:r><w:r><w:t>u</w:t></w:r><w:r><w:t xml:space="preserve"> została podzielona na trzy sekcje</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>'
Replace-ParaXml (31 + $OFF) $p31

# 7) Merge the three "przetestowanie..." runs into one
#    (old index 34 -> now 36).
$p34 = '<w:p><w:pPr><w:pStyle w:val="Bezodstpw"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>przetestowanie aplikacji pod kątem zabezpieczeń, występowania błędów</w:t></w:r></w:p>'
Replace-ParaXml (34 + $OFF) $p34

# 8) Remove the _GoBack bookmark from the "5. Podsumowanie" paragraph
#    (old index 45 -> now 47).
$p45 = '<w:p><w:pPr><w:ind w:firstLine="709"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr><w:t>Udało się zrealizować wszystkie wcześniej ustalone założenia. Największy problem dotyczył skalowania pełnego ekranu na monitorach o różnej rozdzielczości, lecz dzięki współpracy problem udało się rozwiązać i zrealizować projekt do końca. W dalszej perspektywie do aplikacji może zostać dodany tryb gracz vs bot, oraz różne rodzaje bomb.</w:t></w:r></w:p>'
Replace-ParaXml (45 + $OFF) $p45

Write-Output "all edits applied"
